# Generate Report for Handback
# -----------------------------------------------------------------------
# This script reproduces, via Excel COM automation, the "handback" report
# generation step: the localization status workbook is updated so that:
#   - the Overview / per-language "Status" column shows that the files
#     have been handed back and are in sync with en-US
#   - the per-language sheets (zh-cn, de-de) gain a "Latest Target File"
#     hyperlink (column I) and "Latest Handback File" value (column J)
#     for each source file row
#   - the "Latest Handback DateTime" column (K) is stamped with the
#     real handback timestamps
#   - columns that now hold longer text are widened to fit
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

$hyperlinkColor = 15570276   # OLE BGR for RGB(100,149,237) / #6495ED (cornflower blue)

# -----------------------------------------------------------------------
# 1. Status column updates (all cells that used to show "Ready for
#    handoff" now show the handback status). Overview shows one status
#    per language per row; the language sheets show a single Status
#    column.
# -----------------------------------------------------------------------
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("C3").Value = $statusText

$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("C3").Value = $statusText

# -----------------------------------------------------------------------
# 2. zh-cn sheet: Latest Target File (I) + Latest Handback File (J) +
#    Latest Handback DateTime (K)
# -----------------------------------------------------------------------

# Row 2 -> 5cda55df-95e1-436c-b019-5ed0b32e57b5.md
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4c14b2cfdf5859206d8a154898e89587ebc6d078/e2e/5cda55df-95e1-436c-b019-5ed0b32e57b5.md", "", "", "5cda55df-95e1-436c-b019-5ed0b32e57b5.md")
$wsZhCn.Range("I2").Font.Underline = $true
$wsZhCn.Range("I2").Font.Color = $hyperlinkColor
$wsZhCn.Range("J2").Value = "5cda55df-95e1-436c-b019-5ed0b32e57b5.8f26b7bd62c363764a0e1f4eb1f5c79889efba3e.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-19 02:58:07"

# Row 3 -> d1724f0c-1c19-4ace-ad47-e6b615a77776.md
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4c14b2cfdf5859206d8a154898e89587ebc6d078/e2e/d1724f0c-1c19-4ace-ad47-e6b615a77776.md", "", "", "d1724f0c-1c19-4ace-ad47-e6b615a77776.md")
$wsZhCn.Range("I3").Font.Underline = $true
$wsZhCn.Range("I3").Font.Color = $hyperlinkColor
$wsZhCn.Range("J3").Value = "d1724f0c-1c19-4ace-ad47-e6b615a77776.99a40a0743b1ea5f069bea72f4fd569fdc670669.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-19 02:58:07"

# -----------------------------------------------------------------------
# 3. de-de sheet: Latest Target File (I) + Latest Handback File (J) +
#    Latest Handback DateTime (K)
# -----------------------------------------------------------------------

# Row 2 -> 5cda55df-95e1-436c-b019-5ed0b32e57b5.md
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4c14b2cfdf5859206d8a154898e89587ebc6d078/e2e/5cda55df-95e1-436c-b019-5ed0b32e57b5.md", "", "", "5cda55df-95e1-436c-b019-5ed0b32e57b5.md")
$wsDeDe.Range("I2").Font.Underline = $true
$wsDeDe.Range("I2").Font.Color = $hyperlinkColor
$wsDeDe.Range("J2").Value = "5cda55df-95e1-436c-b019-5ed0b32e57b5.8f26b7bd62c363764a0e1f4eb1f5c79889efba3e.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-19 02:58:15"

# Row 3 -> d1724f0c-1c19-4ace-ad47-e6b615a77776.md
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4c14b2cfdf5859206d8a154898e89587ebc6d078/e2e/d1724f0c-1c19-4ace-ad47-e6b615a77776.md", "", "", "d1724f0c-1c19-4ace-ad47-e6b615a77776.md")
$wsDeDe.Range("I3").Font.Underline = $true
$wsDeDe.Range("I3").Font.Color = $hyperlinkColor
$wsDeDe.Range("J3").Value = "d1724f0c-1c19-4ace-ad47-e6b615a77776.99a40a0743b1ea5f069bea72f4fd569fdc670669.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-19 02:58:15"

# -----------------------------------------------------------------------
# 4. Column widths: widen the columns that now hold the longer status /
#    file-name / hyperlink text so the content is not clipped.
#    (ColumnWidth values below are chosen so the saved width lands as
#    close as possible to the intended "fit" width on this engine's
#    width quantization grid.)
# -----------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.15   # E: zh-cn status
$wsOverview.Columns.Item(6).ColumnWidth = 29.15   # F: de-de status

$wsZhCn.Columns.Item(3).ColumnWidth = 29.15        # C: Status
$wsZhCn.Columns.Item(9).ColumnWidth = 39.17         # I: Latest Target File
$wsZhCn.Columns.Item(10).ColumnWidth = 39.17        # J: Latest Handback File

$wsDeDe.Columns.Item(3).ColumnWidth = 29.15        # C: Status
$wsDeDe.Columns.Item(9).ColumnWidth = 39.17         # I: Latest Target File
$wsDeDe.Columns.Item(10).ColumnWidth = 39.17        # J: Latest Handback File
